# Refresh the crypto price/volume snapshot (Price = column D, Volume(1h) = column E)
# for the rows whose source data moved since the last run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.562.75"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "2.458.81"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'318.50"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").Value = "'91.14"
$ws.Range("E6").Value = "  -1.65%  "
$ws.Range("E7").Value = "  -0.66%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.503"
$ws.Range("E9").Value = "  -2.24%  "
$ws.Range("D10").Value = "'0.0850"
$ws.Range("E10").Value = "  -5.63%  "
$ws.Range("D11").Value = "'32.57"
$ws.Range("E11").Value = "  -0.76%  "
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("D13").Value = "2.837.93"
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("E14").Value = "  -1.01%  "
$ws.Range("D15").Value = "'15.40"
$ws.Range("E15").Value = "  -2.01%  "
$ws.Range("D16").Value = "2.525.58"
$ws.Range("E16").Value = "  +2.16%  "
$ws.Range("D17").Value = "'0.783"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "41.497.79"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").Value = "'6.40"
$ws.Range("E19").Value = "  -1.62%  "
$ws.Range("D20").Value = "0.0₃0936"
$ws.Range("E20").Value = "  -4.73%  "
$ws.Range("D21").Value = "'71.95"
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("E22").Value = "  -2.48%  "
$ws.Range("D23").Value = "'237.42"
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("D24").Value = "'2.73"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "'24.59"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("E28").Value = "  -1.67%  "
$ws.Range("D29").Value = "'9.65"
$ws.Range("E29").Value = "  -1.90%  "
$ws.Range("D30").Value = "'36.02"
$ws.Range("E30").Value = "  +2.14%  "
$ws.Range("D31").Value = "'156.59"
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("D32").Value = "'5.39"
$ws.Range("E32").Value = "  -2.02%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("D35").Value = "'0.0760"
$ws.Range("E35").Value = "  -1.11%  "
$ws.Range("D36").Value = "'16.83"
$ws.Range("E36").Value = "  -3.79%  "
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("D38").Value = "'0.115"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  +0.49%  "
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").Value = "'2.32"
$ws.Range("E42").Value = "  -7.65%  "
$ws.Range("D43").Value = "1.998.95"
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("D44").Value = "'0.0281"
$ws.Range("E44").Value = "  -1.39%  "
$ws.Range("D45").Value = "'18.46"
$ws.Range("E45").Value = "  -2.45%  "
$ws.Range("D46").Value = "'2.93"
$ws.Range("E46").Value = "  -0.85%  "
$ws.Range("D47").Value = "'9.51"
$ws.Range("E47").Value = "  +4.73%  "
$ws.Range("D48").Value = "2.718.92"
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("D49").Value = "'75.73"
$ws.Range("E49").Value = "  +4.17%  "
$ws.Range("D50").Value = "'96.72"
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("D51").Value = "'66.43"
$ws.Range("E51").Value = "  -0.78%  "
